# Generate Report for Handback
# The previous handback round had a stale/mismatched-version error; this
# run's report reflects a clean handback that is in sync with en-US, with
# refreshed handback timestamps and the error detail cleared.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # Overview
$ws2 = $wb.Worksheets.Item(2)   # zh-cn
$ws3 = $wb.Worksheets.Item(3)   # de-de

$newStatus = "Handed back: in sync with en-US"

# --- Overview sheet: Status shown per-locale in columns E (zh-cn) and F (de-de) ---
$ws1.Range("E2").Value = $newStatus
$ws1.Range("F2").Value = $newStatus

# --- zh-cn detail sheet ---
$ws2.Range("C2").Value = $newStatus                    # Status
$ws2.Range("L2").Value = "2017-02-15 06:04:42"          # Latest Handback DateTime
$ws2.Range("R2").Value = ""                             # Error Detail cleared

# --- de-de detail sheet ---
$ws3.Range("C2").Value = $newStatus                     # Status
$ws3.Range("L2").Value = "2017-02-15 06:05:09"           # Latest Handback DateTime
$ws3.Range("R2").Value = ""                              # Error Detail cleared

# --- Column widths widened/narrowed to fit the refreshed content ---
$ws1.Columns.Item(5).ColumnWidth = 29.166666666666668    # Overview col E (zh-cn status)
$ws1.Columns.Item(6).ColumnWidth = 29.166666666666668    # Overview col F (de-de status)

$ws2.Columns.Item(3).ColumnWidth = 29.166666666666668    # zh-cn col C (Status)
$ws2.Columns.Item(18).ColumnWidth = 12.833333333333332   # zh-cn col R (Error Detail)

$ws3.Columns.Item(3).ColumnWidth = 29.166666666666668    # de-de col C (Status)
$ws3.Columns.Item(18).ColumnWidth = 12.833333333333332   # de-de col R (Error Detail)
